# Applies the table-style change described by the commit diff:
# three tables (on slides 14, 15 and 16) switch from the deck's local
# custom table style ("Table_0", {551866AD-A619-40DD-9D89-7A4AFBEBFE55})
# to the built-in PowerPoint table style {C36492E9-721D-4A1A-8338-C9732DC83061}.

$p = $ppt.ActivePresentation

$newStyleId = "{C36492E9-721D-4A1A-8338-C9732DC83061}"
$slideIndexesWithTables = @(14, 15, 16)

foreach ($slideIndex in $slideIndexesWithTables) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
